# Applies the "Updated symbol list" data refresh (Tue Jan  3 14:12:13 UTC 2023)
# to the cryptos worksheet: refreshed Price/Volume(1h) figures, the Hora (hour)
# column advancing from 13 -> 14 for every data row, and a reordering of the
# CEJI / KickToken rows (42 and 43 swap Coin name + Link).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row column -> new value map, keyed by row number (matches sheet rows 2-51).
$rowUpdates = [ordered]@{
    2 = @{ 'D'='245.64'; 'E'='-0.40%'; 'G'='14' }
    3 = @{ 'E'='-0.53%'; 'G'='14' }
    4 = @{ 'D'='5.280'; 'E'='2.04%'; 'G'='14' }
    5 = @{ 'D'='0.05750'; 'E'='0.72%'; 'G'='14' }
    6 = @{ 'D'='6.642'; 'E'='0.55%'; 'G'='14' }
    7 = @{ 'D'='3.246'; 'E'='6.27%'; 'G'='14' }
    8 = @{ 'D'='0.8589'; 'E'='0.02%'; 'G'='14' }
    9 = @{ 'D'='0.8535'; 'E'='-2.12%'; 'G'='14' }
    10 = @{ 'D'='0.1379'; 'E'='1.26%'; 'G'='14' }
    11 = @{ 'E'='0.27%'; 'G'='14' }
    12 = @{ 'D'='0.03146'; 'E'='9.92%'; 'G'='14' }
    13 = @{ 'D'='0.09344'; 'E'='-0.45%'; 'G'='14' }
    14 = @{ 'D'='0.001533'; 'E'='0.15%'; 'G'='14' }
    15 = @{ 'D'='0.0005978'; 'E'='-0.04%'; 'G'='14' }
    16 = @{ 'D'='0.005961'; 'E'='-3.47%'; 'G'='14' }
    17 = @{ 'D'='3.517'; 'E'='1.08%'; 'G'='14' }
    18 = @{ 'E'='0.32%'; 'G'='14' }
    19 = @{ 'D'='0.3155'; 'E'='0.35%'; 'G'='14' }
    20 = @{ 'D'='0.03317'; 'E'='1.41%'; 'G'='14' }
    21 = @{ 'E'='0.60%'; 'G'='14' }
    22 = @{ 'D'='3.477'; 'E'='0.16%'; 'G'='14' }
    23 = @{ 'D'='0.1410'; 'E'='2.09%'; 'G'='14' }
    24 = @{ 'D'='0.04109'; 'E'='-1.37%'; 'G'='14' }
    25 = @{ 'D'='0.001221'; 'E'='0.12%'; 'G'='14' }
    26 = @{ 'D'='0.004168'; 'E'='-18.16%'; 'G'='14' }
    27 = @{ 'D'='0.0001199'; 'E'='-0.92%'; 'G'='14' }
    28 = @{ 'D'='0.0001449'; 'E'='-25.28%'; 'G'='14' }
    29 = @{ 'G'='14' }
    30 = @{ 'G'='14' }
    31 = @{ 'G'='14' }
    32 = @{ 'G'='14' }
    33 = @{ 'G'='14' }
    34 = @{ 'G'='14' }
    35 = @{ 'G'='14' }
    36 = @{ 'G'='14' }
    37 = @{ 'G'='14' }
    38 = @{ 'G'='14' }
    39 = @{ 'G'='14' }
    40 = @{ 'D'='0.03775'; 'E'='0.29%'; 'G'='14' }
    41 = @{ 'D'='0.1071'; 'G'='14' }
    42 = @{ 'B'='KickToken'; 'C'='https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'; 'D'='0.003586'; 'E'='-36.77%'; 'G'='14' }
    43 = @{ 'B'='CEJI'; 'C'='https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'; 'D'='0.002393'; 'E'='-5.05%'; 'G'='14' }
    44 = @{ 'D'='0.009216'; 'E'='-6.00%'; 'G'='14' }
    45 = @{ 'D'='0.00005418'; 'E'='6.11%'; 'G'='14' }
    46 = @{ 'E'='-0.03%'; 'G'='14' }
    47 = @{ 'D'='0.08986'; 'E'='26.56%'; 'G'='14' }
    48 = @{ 'D'='0.002828'; 'E'='4.13%'; 'G'='14' }
    49 = @{ 'E'='-0.03%'; 'G'='14' }
    50 = @{ 'E'='-0.03%'; 'G'='14' }
    51 = @{ 'G'='14' }
}

foreach ($row in $rowUpdates.Keys) {
    $cols = $rowUpdates[$row]
    foreach ($col in $cols.Keys) {
        $cell = $ws.Range("$col$row")
        # Force a Text number format so numeric-looking values (prices,
        # percentages, and the Hora column) are stored verbatim instead of
        # being auto-converted to floating point numbers by Excel.
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
    }
}

